# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worker data table (B15:J18) gets re-sorted: the existing two rows for
# JAVIER EDUARDO PUERTA COLINA (periods 2402 and 2401) are reflowed so that
# the 2401 period now appears first (row 16), the DOVANIS DE JESUS MONTAÑO
# VIADERO / CC / period 2401 row moves up to row 17, and the JAVIER /
# period 2402 row moves down to become the last data row (row 18).
#
# Net effect per row (columns B..G only change; H/I/J stay blank):
#   Row 16: PPT | 7448501    | JAVIER EDUARDO PUERTA COLINA      | 2401 | 48000 | 1200000
#   Row 17: CC  | 1127591543 | DOVANIS DE JESUS MONTAÑO VIADERO  | 2401 | 6933  | 1300000
#   Row 18: PPT | 7448501    | JAVIER EDUARDO PUERTA COLINA      | 2402 | 1600  | 1200000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> PPT / 7448501 / JAVIER.../ 2401 / 48000 / 1200000
$ws.Range("B16").Value = "PPT"
$ws.Range("C16").Value = "7448501"
$ws.Range("D16").Value = "JAVIER EDUARDO PUERTA COLINA"
$ws.Range("E16").Value = "2401"
$ws.Range("F16").Value = 48000
$ws.Range("G16").Value = 1200000

# Row 17 -> CC / 1127591543 / DOVANIS.../ 2401 / 6933 / 1300000
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1127591543"
$ws.Range("D17").Value = "DOVANIS DE JESUS MONTAÑO VIADERO"
$ws.Range("E17").Value = "2401"
$ws.Range("F17").Value = 6933
$ws.Range("G17").Value = 1300000

# Row 18 -> PPT / 7448501 / JAVIER.../ 2402 / 1600 / 1200000
$ws.Range("B18").Value = "PPT"
$ws.Range("C18").Value = "7448501"
$ws.Range("D18").Value = "JAVIER EDUARDO PUERTA COLINA"
$ws.Range("E18").Value = "2402"
$ws.Range("F18").Value = 1600
$ws.Range("G18").Value = 1200000

$wb.Save()
